$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.198.71"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "2.642.87"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.26"
$ws.Range("D5").Style = $ws.Range("C2").Style
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.45"
$ws.Range("D6").Style = $ws.Range("C2").Style
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = $ws.Range("C2").Style
$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "2.649.74"
$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.91"
$ws.Range("D10").Style = $ws.Range("C2").Style
$ws.Range("E10").Value = "  +9.89%  "

$ws.Range("E11").Value = "  -2.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("D12").Style = $ws.Range("C2").Style
$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("E13").Value = "  +2.04%  "

$ws.Range("D14").Value = "3.110.50"
$ws.Range("E14").Value = "  +0.08%  "

$ws.Range("D15").Value = "59.249.66"
$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.15"
$ws.Range("D16").Style = $ws.Range("C2").Style
$ws.Range("E16").Value = "  +0.92%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.668.02"
$ws.Range("E17").Value = "  +0.55%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000136"
$ws.Range("D18").Style = $ws.Range("C2").Style
$ws.Range("E18").Value = "  -1.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "341.77"
$ws.Range("D19").Style = $ws.Range("C2").Style
$ws.Range("E19").Value = "  -2.34%  "

$ws.Range("E20").Value = "  -1.53%  "

$ws.Range("E21").Value = "  +0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.30"
$ws.Range("D22").Style = $ws.Range("C2").Style
$ws.Range("E22").Value = "  +1.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = $ws.Range("C2").Style
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.71"
$ws.Range("D24").Style = $ws.Range("C2").Style
$ws.Range("E24").Value = "  +2.66%  "

$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("E26").Value = "  -1.02%  "

$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("D28").Value = "0.0₃0802"
$ws.Range("E28").Value = "  -0.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.11"
$ws.Range("D29").Style = $ws.Range("C2").Style
$ws.Range("E29").Value = "  -0.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.67"
$ws.Range("D30").Style = $ws.Range("C2").Style
$ws.Range("E30").Value = "  +2.66%  "

$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("E32").Value = "  +0.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.78"
$ws.Range("D33").Style = $ws.Range("C2").Style
$ws.Range("E33").Value = "  -0.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.32"
$ws.Range("D34").Style = $ws.Range("C2").Style
$ws.Range("E34").Value = "  -0.30%  "

$ws.Range("E35").Value = "  +3.13%  "

$ws.Range("E36").Value = "  +1.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.892"
$ws.Range("D37").Style = $ws.Range("C2").Style
$ws.Range("E37").Value = "  -6.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.884"
$ws.Range("D38").Style = $ws.Range("C2").Style
$ws.Range("E38").Value = "  +2.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.80"
$ws.Range("D39").Style = $ws.Range("C2").Style
$ws.Range("E39").Value = "  +0.47%  "

$ws.Range("E40").Value = "  +0.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.59"
$ws.Range("D41").Style = $ws.Range("C2").Style
$ws.Range("E41").Value = "  -2.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.626"
$ws.Range("D42").Style = $ws.Range("C2").Style
$ws.Range("E42").Value = "  +4.22%  "

$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "275.71"
$ws.Range("D44").Style = $ws.Range("C2").Style
$ws.Range("E44").Value = "  -0.65%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.81"
$ws.Range("D45").Style = $ws.Range("C2").Style
$ws.Range("E45").Value = "  +0.87%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0977"
$ws.Range("D46").Style = $ws.Range("C2").Style
$ws.Range("E46").Value = "  -1.74%  "

$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0540"
$ws.Range("D47").Style = $ws.Range("C2").Style
$ws.Range("E47").Value = "  +1.96%  "

$ws.Range("D48").Value = "2.051.03"
$ws.Range("E48").Value = "  -1.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.52"
$ws.Range("D49").Style = $ws.Range("C2").Style
$ws.Range("E49").Value = "  +1.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.80"
$ws.Range("D50").Style = $ws.Range("C2").Style
$ws.Range("E50").Value = "  +1.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.97"
$ws.Range("D51").Style = $ws.Range("C2").Style
$ws.Range("E51").Value = "  -0.36%  "
